$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: Monday 4-Sep-2023 ---
$ws.Range("A16").Value = 45173
$ws.Range("A16").NumberFormat = "d-mmm"
$ws.Range("B16:E16").Value = "PRESENT"
$ws.Range("F16:I16").Value = "ABSENT"

$ws.Range("F16").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("G16").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("H16").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("I16").AddComment("RENUKA:`nNo Response") | Out-Null

# --- Row 17: Tuesday 5-Sep-2023 ---
$ws.Range("A17").Value = 45174
$ws.Range("A17").NumberFormat = "d-mmm"
$ws.Range("B17:E17").Value = "PRESENT"
$ws.Range("F17:G17").Value = "ABSENT"
$ws.Range("H17").Value = "PRESENT"
$ws.Range("I17").Value = "ABSENT"

$ws.Range("F17").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("G17").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("I17").AddComment("RENUKA:`nNo Response") | Out-Null

# --- View state: scroll so column C is leftmost, selection on I17 ---
$ws.Range("I17").Select()
